# "Generate Report for Handback"
#
# Refresh the handoff/handback timestamp columns on the per-language
# report sheets (these are the "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" columns, E and H respectively) to
# reflect a freshly (re-)generated report.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: handoff 17:20:27 -> 17:22:07, handback 17:21:14 -> 17:22:48
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-19 17:22:07"
$wsZh.Range("E3").Value = "2016-03-19 17:22:07"
$wsZh.Range("H2").Value = "2016-03-19 17:22:48"
$wsZh.Range("H3").Value = "2016-03-19 17:22:48"

# de-de sheet: handoff 17:20:43 -> 17:22:15, handback 17:21:27 -> 17:23:01
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-19 17:22:15"
$wsDe.Range("E3").Value = "2016-03-19 17:22:15"
$wsDe.Range("H2").Value = "2016-03-19 17:23:01"
$wsDe.Range("H3").Value = "2016-03-19 17:23:01"
